$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, bordered, centered) from N1 onto the new header cells O1:P1
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column O = "t+15" values, rows 1-67 (row 1 is the header "14")
$oValues = @(
    14,
    -0.2128477638078029,
    0.2813075268075425,
    0.2463629988121556,
    -0.0307656133348515,
    0.2306213330448918,
    -0.371172054729902,
    -0.1429756724898395,
    -0.2833047036358718,
    0.4306152552986454,
    -0.1779404439221652,
    -0.00554741234322996,
    -0.01994034920993558,
    0.2096083778699824,
    0.1158366582466219,
    0.5667400886097449,
    0.6435036435735713,
    -0.06619109519494865,
    0.4147082358362778,
    0.4680111208729821,
    0.6795640375324163,
    0.4491915327500307,
    -0.005571332609762891,
    2.157780851749858,
    0.2917440745863075,
    0.152467296036162,
    0.04714036940376232,
    0.7963412485181309,
    1.954389059321646,
    0.640684197316871,
    -0.4924788820920559,
    0.5400243399792264,
    0.7445387718341611,
    -0.8560527810546902,
    0.7731180761019696,
    0.710445870074472,
    0.6734072349369168,
    0.6235520964690628,
    0.5754198577815065,
    0.7384154861031726,
    0.5485594665135802,
    0.5677140889241663,
    0.6492034822425543,
    0.6634945156418861,
    0.6184362577945474,
    -1.292567223126907,
    -1.009902242279011,
    -0.884853480017777,
    -0.6459783835648634,
    -0.0628550302052869,
    -0.8724581358318191,
    -0.8724581358318191,
    -1.147104558077292,
    -0.1780386108715763,
    -1.029059925247228,
    -0.8968228105385498,
    -0.9031757676071503,
    -1.047694685311832,
    -0.7824762314578961,
    -0.403574102318818,
    0.3874036144745389,
    -1.153390528667371,
    -0.5641777544994385,
    -0.8394084087706666,
    -0.02225654199021455,
    -0.7305266003811319,
    -0.6984980672173858
)

# Column P = "t+16" values, rows 1-67 (row 1 is the header "15")
$pValues = @(
    15,
    -0.2118901065629061,
    0.2811752558168016,
    0.2486797524839374,
    -0.03225841381480378,
    0.2299902910042842,
    -0.3699383689509853,
    -0.1384853716517394,
    -0.2800782769518634,
    0.4294587097203348,
    -0.1776001658684229,
    -0.01069848376877072,
    -0.02223238652636863,
    0.2011992897294233,
    0.1019751367546734,
    0.5523092529592173,
    0.626826037000103,
    -0.07453402987895902,
    0.4059315240000756,
    0.4512426545136669,
    0.6662519438824651,
    0.4349030728676286,
    -0.01805577677420292,
    1.959933607157823,
    0.2896687705687532,
    0.1441163917052351,
    0.03805970042251687,
    0.7886533795173142,
    1.819249605817164,
    0.6334882871403437,
    -0.4926801269012133,
    0.5336091816930694,
    0.7413142241910726,
    -0.8599548809714945,
    0.7750422638657158,
    0.7148727103196514,
    0.6778230668963577,
    0.6237819527619272,
    0.5778518553257703,
    0.7404350324146238,
    0.5515814593899226,
    0.5704519990425972,
    0.6507580129854277,
    0.6670905577338264,
    0.6257897093034609,
    -1.29482394536213,
    -1.01158835404017,
    -0.8843913445364918,
    -0.6458179574229141,
    -0.06427546707599818,
    -0.8713647033370471,
    -0.8713647033370471,
    -1.147667061955454,
    -0.1781667638844488,
    -1.031078113966504,
    -0.8999024290092584,
    -0.9088797600531041,
    -1.052485398801365,
    -0.7830473127192779,
    -0.4052453470076435,
    0.3872804548219493,
    -1.159886010957426,
    -0.5602152825374364,
    -0.8396217037502596,
    -0.02334129164106661,
    -0.7371199583381679,
    -0.7075701313924574
)

for ($i = 0; $i -lt $oValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 15).Value = $oValues[$i]
    $ws.Cells.Item($i + 1, 16).Value = $pValues[$i]
}